$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the " IT DEPARTMENT" value from A4 up to A2, leaving A4 blank.
$deptText = $ws.Range("A4").Text
$ws.Range("A2").Value = $deptText
$ws.Range("A2").Style = "Normal 2"
$ws.Range("A4").ClearContents()

# Remove the now-extra trailing rows (7, 8 and 9) entirely.
$ws.Rows("7:9").Delete()

# Restore the original selection to A2.
$ws.Range("A2").Select()
